# 1st changes of mifos to finflux
# Insert a new (blank) column before column N ("Late") on the
# "Repayment schedule" worksheet, pushing the existing N/O/P columns
# (Late / Original / Outstanding) one place to the right (-> O/P/Q).
# The new column inherits the column width of its left neighbour (M).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new column at N; existing N (and everything to the right)
# shifts right to O, O->P, P->Q.
$ws.Columns("N").Insert()

# Give the freshly inserted column the same width as column M, just like
# Excel does when a column is inserted (format copied from the column to
# the left).
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

# Leave the final selection on cell M13, matching the saved selection
# state recorded for this sheet.
[void]$ws.Range("M13").Select()
